$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("<made>",    "<made>",    18),
    @("<about>",   "<about>",   13),
    @("<have>",    "<have>",    12),
    @("<an>",      "<an>",      14),
    @("<nine>",    "<nine>",    18),
    @("<left>",    "<left>",    12),
    @("<oscar>",   "<oscar>",   17),
    @("<line>",    "<by>",      22),
    @("<alt>",     "<alt>",     8),
    @("<see>",     "<see>",     15),
    @("<seven>",   "<seven>",   13),
    @("<command>", "<command>", 15),
    @("<day>",     "<they>",    14),
    @("<victor>",  "<victor>",  15),
    @("<we>",      "<we>",      13),
    @("<for>",     "<four>",    15),
    @("<out>",     "<how>",     13)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
